$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The BOM/CPL had a stale duplicate line for the analog switch: row 7 listed
# the old "ADG702 SPST-NC switch" part (U2, C133478) while row 8 already
# carried the correct "TS5A3160 SPDT Analog Switch" part (U2, C185770) for
# the same designator/footprint. Remove the obsolete row so the component
# orientation/part data lines up correctly for the CPL export; the rows
# below (U3 / op-amp) shift up to fill the gap.
$ws.Rows.Item(7).Delete()

# Restore the last sheet row's cached row-height metadata: deleting a row
# shifts every row below it up by one, and on this sheet every row all the
# way down to row 1048576 carried an explicit (non-default) row height.
$ws.Rows.Item(1048576).RowHeight = 12.8

# Leave the selection where the saved workbook had it.
$ws.Range("A7").Select()
